$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime (D2) and Correspond Handback DateTime (G2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-26 09:08:38"
$wsZhCn.Range("G2").Value = "2016-01-26 09:09:22"

# de-de sheet: update Correspond Handoff Datetime (D2) and Correspond Handback DateTime (G2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-26 09:08:49"
$wsDeDe.Range("G2").Value = "2016-01-26 09:09:40"
